$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.08758366666666667
$ws.Range("H2").Value = 0.262751
$ws.Range("I2").Value = 0.2371976925785164
$ws.Range("J2").Value = 0.2371976925785164
$ws.Range("M2").Value = 42.105049
$ws.Range("N2").Value = 126.315147
$ws.Range("O2").Value = 0.596182887750605
$ws.Range("P2").Value = 0.5961828877506051
$ws.Range("Q2").Value = 3.687714576599667
$ws.Range("R2").Value = 33.18943118939701
$ws.Range("S2").Value = 0.1414132053292402
$ws.Range("T2").Value = 0.1414132053292402
$ws.Range("G3").Value = 0.08758366666666667
$ws.Range("H3").Value = 0.262751
$ws.Range("I3").Value = 0.2371976925785164
$ws.Range("J3").Value = 0.2371976925785164
$ws.Range("O3").Value = 0.04140484982922635
$ws.Range("P3").Value = 0.04140484982922635
$ws.Range("Q3").Value = 0.2561114573973333
$ws.Range("R3").Value = 2.305003116576
$ws.Range("S3").Value = 0.009821134841052469
$ws.Range("T3").Value = 0.009821134841052469
$ws.Range("G4").Value = 0.08758366666666667
$ws.Range("H4").Value = 0.262751
$ws.Range("I4").Value = 0.2371976925785164
$ws.Range("J4").Value = 0.2371976925785164
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.703340666666667
$ws.Range("N4").Value = 23.110022
$ws.Range("O4").Value = 0.1090748020262369
$ws.Range("P4").Value = 0.1090748020262369
$ws.Range("Q4").Value = 0.6746868211691112
$ws.Range("R4").Value = 6.072181390522
$ws.Range("S4").Value = 0.02587229135908188
$ws.Range("T4").Value = 0.02587229135908188
$ws.Range("G5").Value = 0.08758366666666667
$ws.Range("H5").Value = 0.262751
$ws.Range("I5").Value = 0.2371976925785164
$ws.Range("J5").Value = 0.2371976925785164
$ws.Range("M5").Value = 16.80268266666667
$ws.Range("N5").Value = 50.40804800000001
$ws.Range("O5").Value = 0.2379161671126513
$ws.Range("P5").Value = 0.2379161671126513
$ws.Range("Q5").Value = 1.471640557783111
$ws.Range("R5").Value = 13.244765020048
$ws.Range("S5").Value = 0.0564331658662456
$ws.Range("T5").Value = 0.0564331658662456
$ws.Range("G6").Value = 0.08758366666666667
$ws.Range("H6").Value = 0.262751
$ws.Range("I6").Value = 0.2371976925785164
$ws.Range("J6").Value = 0.2371976925785164
$ws.Range("M6").Value = 0.6492376666666667
$ws.Range("N6").Value = 1.947713
$ws.Range("O6").Value = 0.009192825947068677
$ws.Range("P6").Value = 0.009192825947068679
$ws.Range("Q6").Value = 0.05686261538477778
$ws.Range("R6").Value = 0.511763538463
$ws.Range("S6").Value = 0.002180517102920605
$ws.Range("T6").Value = 0.002180517102920606
$ws.Range("G7").Value = 0.08758366666666667
$ws.Range("H7").Value = 0.262751
$ws.Range("I7").Value = 0.2371976925785164
$ws.Range("J7").Value = 0.2371976925785164
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4398816666666667
$ws.Range("N7").Value = 1.319645
$ws.Range("O7").Value = 0.006228467334211686
$ws.Range("P7").Value = 0.006228467334211686
$ws.Range("Q7").Value = 0.03852644926611112
$ws.Range("R7").Value = 0.346738043395
$ws.Range("S7").Value = 0.001477378079975675
$ws.Range("T7").Value = 0.001477378079975675
$ws.Range("I8").Value = 0.4417069141397272
$ws.Range("J8").Value = 0.4417069141397272
$ws.Range("M8").Value = 42.105049
$ws.Range("N8").Value = 126.315147
$ws.Range("O8").Value = 0.596182887750605
$ws.Range("P8").Value = 0.5961828877506051
$ws.Range("Q8").Value = 6.867221211769334
$ws.Range("R8").Value = 61.80499090592401
$ws.Range("S8").Value = 0.2633381036112311
$ws.Range("T8").Value = 0.2633381036112312
$ws.Range("I9").Value = 0.4417069141397272
$ws.Range("J9").Value = 0.4417069141397272
$ws.Range("O9").Value = 0.04140484982922635
$ws.Range("P9").Value = 0.04140484982922635
$ws.Range("S9").Value = 0.01828880844848638
$ws.Range("T9").Value = 0.01828880844848638
$ws.Range("I10").Value = 0.4417069141397272
$ws.Range("J10").Value = 0.4417069141397272
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.703340666666667
$ws.Range("N10").Value = 23.110022
$ws.Range("O10").Value = 0.1090748020262369
$ws.Range("P10").Value = 0.1090748020262369
$ws.Range("Q10").Value = 1.256394320491556
$ws.Range("R10").Value = 11.307548884424
$ws.Range("S10").Value = 0.04817909421341076
$ws.Range("T10").Value = 0.04817909421341077
$ws.Range("I11").Value = 0.4417069141397272
$ws.Range("J11").Value = 0.4417069141397272
$ws.Range("M11").Value = 16.80268266666667
$ws.Range("N11").Value = 50.40804800000001
$ws.Range("O11").Value = 0.2379161671126513
$ws.Range("P11").Value = 0.2379161671126513
$ws.Range("Q11").Value = 2.740472735779556
$ws.Range("R11").Value = 24.664254622016
$ws.Range("S11").Value = 0.1050892159992809
$ws.Range("T11").Value = 0.1050892159992809
$ws.Range("I12").Value = 0.4417069141397272
$ws.Range("J12").Value = 0.4417069141397272
$ws.Range("M12").Value = 0.6492376666666667
$ws.Range("N12").Value = 1.947713
$ws.Range("O12").Value = 0.009192825947068677
$ws.Range("P12").Value = 0.009192825947068679
$ws.Range("Q12").Value = 0.1058889321328889
$ws.Range("R12").Value = 0.9530003891960002
$ws.Range("S12").Value = 0.00406053478130332
$ws.Range("T12").Value = 0.004060534781303321
$ws.Range("I13").Value = 0.4417069141397272
$ws.Range("J13").Value = 0.4417069141397272
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4398816666666667
$ws.Range("N13").Value = 1.319645
$ws.Range("O13").Value = 0.006228467334211686
$ws.Range("P13").Value = 0.006228467334211686
$ws.Range("Q13").Value = 0.07174352681555556
$ws.Range("R13").Value = 0.64569174134
$ws.Range("S13").Value = 0.002751157086014737
$ws.Range("T13").Value = 0.002751157086014737
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1185623333333333
$ws.Range("H14").Value = 0.355687
$ws.Range("I14").Value = 0.3210953932817563
$ws.Range("J14").Value = 0.3210953932817564
$ws.Range("M14").Value = 42.105049
$ws.Range("N14").Value = 126.315147
$ws.Range("O14").Value = 0.596182887750605
$ws.Range("P14").Value = 0.5961828877506051
$ws.Range("Q14").Value = 4.992072854554333
$ws.Range("R14").Value = 44.92865569098901
$ws.Range("S14").Value = 0.1914315788101337
$ws.Range("T14").Value = 0.1914315788101338
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1185623333333333
$ws.Range("H15").Value = 0.355687
$ws.Range("I15").Value = 0.3210953932817563
$ws.Range("J15").Value = 0.3210953932817564
$ws.Range("O15").Value = 0.04140484982922635
$ws.Range("P15").Value = 0.04140484982922635
$ws.Range("Q15").Value = 0.3466990266346666
$ws.Range("R15").Value = 3.120291239712
$ws.Range("S15").Value = 0.0132949065396875
$ws.Range("T15").Value = 0.0132949065396875
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1185623333333333
$ws.Range("H16").Value = 0.355687
$ws.Range("I16").Value = 0.3210953932817563
$ws.Range("J16").Value = 0.3210953932817564
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 7.703340666666667
$ws.Range("N16").Value = 23.110022
$ws.Range("O16").Value = 0.1090748020262369
$ws.Range("P16").Value = 0.1090748020262369
$ws.Range("Q16").Value = 0.9133260439015556
$ws.Range("R16").Value = 8.219934395114
$ws.Range("S16").Value = 0.03502341645374425
$ws.Range("T16").Value = 0.03502341645374426
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1185623333333333
$ws.Range("H17").Value = 0.355687
$ws.Range("I17").Value = 0.3210953932817563
$ws.Range("J17").Value = 0.3210953932817564
$ws.Range("M17").Value = 16.80268266666667
$ws.Range("N17").Value = 50.40804800000001
$ws.Range("O17").Value = 0.2379161671126513
$ws.Range("P17").Value = 0.2379161671126513
$ws.Range("Q17").Value = 1.992165263219556
$ws.Range("R17").Value = 17.92948736897601
$ws.Range("S17").Value = 0.07639378524712483
$ws.Range("T17").Value = 0.07639378524712485
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.1185623333333333
$ws.Range("H18").Value = 0.355687
$ws.Range("I18").Value = 0.3210953932817563
$ws.Range("J18").Value = 0.3210953932817564
$ws.Range("M18").Value = 0.6492376666666667
$ws.Range("N18").Value = 1.947713
$ws.Range("O18").Value = 0.009192825947068677
$ws.Range("P18").Value = 0.009192825947068679
$ws.Range("Q18").Value = 0.0769751326478889
$ws.Range("R18").Value = 0.6927761938310002
$ws.Range("S18").Value = 0.002951774062844751
$ws.Range("T18").Value = 0.002951774062844752
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.1185623333333333
$ws.Range("H19").Value = 0.355687
$ws.Range("I19").Value = 0.3210953932817563
$ws.Range("J19").Value = 0.3210953932817564
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.4398816666666667
$ws.Range("N19").Value = 1.319645
$ws.Range("O19").Value = 0.006228467334211686
$ws.Range("P19").Value = 0.006228467334211686
$ws.Range("Q19").Value = 0.05215339679055556
$ws.Range("R19").Value = 0.469380571115
$ws.Range("S19").Value = 0.001999932168221274
$ws.Range("T19").Value = 0.001999932168221274
